$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# The "ParagraphCollection.getLast" snippet row had been mapped into the
# table one row too early (row 71, under the "Paragraph" class) instead of
# at the end of the "Paragraph"/"ParagraphCollection" block (row 84), which
# bumped every row from 72-84 up by one slot. Fix: remove the stray row 71
# (the rest of the Paragraph rows shift up to fill the gap, carrying their
# own formatting with them) and re-insert the correct
# ParagraphCollection.getLast row at 84, right after the real last
# "Paragraph" row.

$ws.Rows.Item(71).Delete()
$ws.Rows.Item(84).Insert()

# Deleting/inserting whole rows can leave the worksheet's Excel Table
# (ListObject) one row short of the real used range -- put it back so the
# table still covers A1:F103 like before.
$lo = $ws.ListObjects.Item(1)
if ($lo.Range.Rows.Count -lt 103) {
    $lo.Resize($ws.Range("A1:F103"))
}

$ws.Cells.Item(84, 1).Value2 = "Word"
$ws.Cells.Item(84, 2).Value2 = "ParagraphCollection"
$ws.Cells.Item(84, 3).Value2 = "getLast"
$ws.Cells.Item(84, 4).Value2 = 1
$ws.Cells.Item(84, 5).Value2 = "word-images-insert-and-get-pictures"
$ws.Cells.Item(84, 6).Value2 = "insertImage"

# Column E on this table uses the "vertical-center" cell style; match the
# rest of the block (this mirrors what row 71 originally looked like).
$ws.Cells.Item(84, 5).VerticalAlignment = -4108

# Reflect the editor's final cursor position on the sheet.
$ws.Range("B84").Select() | Out-Null
